$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 79, shifting existing rows (79..207) down to (80..208)
$ws.Rows(79).Insert()

# Populate the newly inserted row 79 with the new data point
$ws.Range("A79").Value = 3
$ws.Range("B79").Value = "Femacal de La Calera"
$ws.Range("C79").Value = "Coquimbo"
$ws.Range("D79").Value = 44645
$ws.Range("E79").Value = 5
$ws.Range("F79").Value = 100112010
$ws.Range("G79").Value = "Achicoria"
$ws.Range("H79").Value = "Sin especificar"
$ws.Range("I79").Value = "Primera"
$ws.Range("J79").Value = 125
$ws.Range("K79").Value = 7500
$ws.Range("L79").Value = 8000
$ws.Range("M79").Value = 7760
$ws.Range("N79").Value = "$/caja 16 unidades"
$ws.Range("O79").Value = "Provincia de Quillota"
$ws.Range("P79").Value = 485
$ws.Range("Q79").Value = 16
$ws.Range("R79").Value = "Hortaliza"
